# Fill in the missing "Generation Charge" (column D) values on the
# "Historical GC" sheet, and correct a handful of existing values that were
# computed from the wrong source row. Values are stored as text (matching
# the rest of the column), so NumberFormat is forced to "@" (Text) before
# assignment to stop Excel from re-interpreting the numeric-looking string
# as a number and dropping its trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historical GC")

$updates = @{
    "D8"   = "8.3750"
    "D9"   = "8.3750"
    "D35"  = "13.1332"
    "D36"  = "13.1332"
    "D47"  = "6.1922"
    "D48"  = "6.1922"
    "D52"  = "7.3430"
    "D53"  = "7.3430"
    "D54"  = "5.8076"
    "D55"  = "5.8076"
    "D66"  = "5.2377"
    "D67"  = "5.2377"
    "D74"  = "5.6646"
    "D75"  = "5.6646"
    "D90"  = "6.2054"
    "D91"  = "6.2054"
    "D123" = "8.2890"
    "D124" = "8.2890"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
